$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166321516036987
$ws.Range("B1").Value = 1.897526860237122
$ws.Range("D1").Value = 2.123118877410889
$ws.Range("E1").Value = 1.066710710525513
